# Correction orthographique pour le poster
#
#  1) "coté" -> "côté" in the "Solutions retenues" paragraph. In the
#     canonical OOXML this shows up as the big run getting split into three
#     runs around the corrected letter.
#  2) The "_GoBack" last-edit bookmark moves from the end of the final
#     paragraph into the middle of the "dans le cadre du projet P2 a
#     l'HE-Arc" run (a side effect of where the cursor was left after the
#     edit), splitting that run into "d" / bookmark / "ans le cadre...".
#
# Bookmarks split a run cleanly at an exact character position without
# touching any other run in the paragraph, whereas writing Range.Text
# re-flows (coalesces) every same-formatted run it touches across the
# whole edited span. So: drop temporary "guard" bookmarks at the run
# boundaries we must not disturb, rewrite the text inside that guarded
# span (so xml:space gets recomputed correctly), then carve the interior
# split point with one more bookmark, and finally remove the scratch
# bookmarks again.

$d = $word.ActiveDocument

# ===== Change 1: paragraph 3 ("Une representation graphique ..."): coté -> côté =====
$p1 = $d.Paragraphs(3)
$p1Start = $p1.Range.Start
$full1 = $p1.Range.Text

$idxBigStart = $full1.IndexOf("résentation")        # start of the big run (after "Une rep")
$posBigStart = $p1Start + $idxBigStart
$idxComplete = $full1.IndexOf("complète")            # start of the next existing run
$posComplete = $p1Start + $idxComplete

# Guard the pre-existing run boundaries on both sides of our edit so the
# text rewrite below cannot bleed into "Une rep" or "complète".
$guard0 = $d.Bookmarks.Add("ZZGUARD0", $d.Range($posBigStart, $posBigStart))
$guard1 = $d.Bookmarks.Add("ZZGUARD1", $d.Range($posComplete, $posComplete))

$idxSplitA = $full1.IndexOf("premier") + 5           # "...dans un premi" | "er temps..."
$posSplitA = $p1Start + $idxSplitA
$idxO = $full1.IndexOf("coté") + 1                   # the "o" that becomes "ô"
$posO = $p1Start + $idxO
$posSplitB = $posO + 1                                # "...cô" | "té..."

# Guard the right edge of the rewritten span too.
$guardB = $d.Bookmarks.Add("ZZSPLITB", $d.Range($posSplitB, $posSplitB))

$leftText = $full1.Substring($idxBigStart, $idxSplitA - $idxBigStart)
$midText  = $full1.Substring($idxSplitA, $idxO - $idxSplitA)

$rAll = $d.Range($posBigStart, $posSplitB)
$rAll.Text = $leftText + $midText + "ô"

# Carve the interior "premi" / "er temps..." boundary.
$guardA = $d.Bookmarks.Add("ZZSPLITA", $d.Range($posSplitA, $posSplitA))

$d.Bookmarks("ZZGUARD0").Delete()
$d.Bookmarks("ZZGUARD1").Delete()
$d.Bookmarks("ZZSPLITA").Delete()
$d.Bookmarks("ZZSPLITB").Delete()

# ===== Change 2: paragraph 5: move "_GoBack" into "dans le cadre du projet P2 à l'HE-Arc" =====
$p2 = $d.Paragraphs(5)
$p2Start = $p2.Range.Start
$full2 = $p2.Range.Text
$idxBm = $full2.IndexOf("dans le cadre") + 1          # right after the leading "d"
$posBm = $p2Start + $idxBm

$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($posBm, $posBm)) | Out-Null
